$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- String constants (existing strategy text reused verbatim, plus new variants) ----
$descRsi = 'content="1️⃣ Explanation: **TSLA Intraday Momentum Reversal Strategy** - Generate ''signal'' based on 5-min bar momentum: Buy if Close > High of 3 periods ago & RSI(10) < 70, Sell if Close < Low of 3 periods ago & RSI(10) > 30, else Hold.

'
$descNew = 'content="1️⃣ Explanation: **TSLA Intraday Momentum Reversal Strategy** - Generates ''signal'' based on 5-min bar momentum: Buy when price closes above the 3-period high after a 2-period downtrend, Sell when below the 3-period low after a 2-period uptrend, else Hold.
'
$codeRsi = '
def add_signal(df):
    df[''CloseShift3''] = df[''Close''].shift(3)
    df[''HighShift3''] = df[''High''].shift(3)
    df[''LowShift3''] = df[''Low''].shift(3)
    df[''DeltaClose''] = df[''Close''] - df[''CloseShift3'']
    df[''RSI_Period''] = 10
    df[''Gain''] = df[''Close''].diff().clip(lower=0)
    df[''Loss''] = -df[''Close''].diff().clip(upper=0)
    df[''AvgGain''] = df[''Gain''].rolling(df[''RSI_Period''].iloc[0]).mean()
    df[''AvgLoss''] = df[''Loss''].rolling(df[''RSI_Period''].iloc[0]).mean()
    df[''RS''] = df[''AvgGain''] / df[''AvgLoss'']
    df[''RSI''] = 100 - (100 / (1 + df[''RS'']))
    buy_condition = (df[''Close''] > df[''HighShift3'']) & (df[''RSI''] < 70)
    sell_condition = (df[''Close''] < df[''LowShift3'']) & (df[''RSI''] > 30)
    df[''signal''] = 0
    df.loc[buy_condition, ''signal''] = 1
    df.loc[sell_condition, ''signal''] = -1
    df.drop(columns=[''CloseShift3'', ''HighShift3'', ''LowShift3'', 
                      ''DeltaClose'', ''RSI_Period'', ''Gain'', ''Loss'', 
                      ''AvgGain'', ''AvgLoss'', ''RS'', ''RSI''], inplace=True)
'
$codeNew1 = '
def add_signal(df):
    df[''Up''] = np.where(df[''Close''] > df[''Close''].shift(1), 1, 0)
    up_days = df[''Up''].rolling(2).sum()
    down_days = df[''Up''].rolling(2).sum().apply(lambda x: 2 - x)
    df[''RecentUptrend''] = np.where(up_days == 2, 1, 0)
    df[''RecentDowntrend''] = np.where(down_days == 2, 1, 0)
    df[''High3''] = df[''High''].rolling(3).max()
    df[''Low3''] = df[''Low''].rolling(3).min()
    df[''signal''] = np.where((df[''Close''] > df[''High3'']) & df[''RecentDowntrend''], 1,
                            np.where((df[''Close''] < df[''Low3'']) & df[''RecentUptrend''], -1, 0))
'
$codeNew2 = '
def add_signal(df):
    df[''pct_change''] = df[''Close''].pct_change().fillna(0)
    up_days = df[''pct_change''] > 0
    down_days = ~up_days
    df[''recent_uptrend''] = (up_days.shift(1)) & (up_days.shift(2))
    df[''recent_downtrend''] = (down_days.shift(1)) & (down_days.shift(2))
    df[''above_3high''] = df[''Close''] > df[''High''].shift(1).rolling(3).max()
    df[''below_3low''] = df[''Close''] < df[''Low''].shift(1).rolling(3).min()
    df[''signal''] = -1
    df.loc[(df[''recent_downtrend'']) & (df[''above_3high'']), ''signal''] = 1
    df.loc[(df[''recent_uptrend'']) & (df[''below_3low'']), ''signal''] = -1
    df.loc[df[''signal''] == -1, ''signal''] = 0
'
$mockImproved = '# mock improved code
def add_signal(df):
    df[''signal'']=0
    return df'

# ---- Update existing row 2 (timestamp + orig_code/improved_code references only change) ----
$ws.Cells.Item(2, 1).Value = 45791.58182560185
$ws.Cells.Item(2, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 3).Value = $codeRsi
$ws.Cells.Item(2, 5).Value = $mockImproved
$ws.Rows.Item(2).AutoFit()

# ---- New rows 3-6 (full new log entries) ----
# Row 3
$ws.Cells.Item(3, 1).Value = 45791.59323685186
$ws.Cells.Item(3, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3, 2).Value = $descRsi
$ws.Cells.Item(3, 3).Value = $codeRsi
$ws.Cells.Item(3, 5).Value = $mockImproved
$ws.Cells.Item(3, 6).Value = 10
$ws.Cells.Item(3, 7).Value = 10
$ws.Cells.Item(3, 8).Value = 2.73
$ws.Cells.Item(3, 9).Value = 49.95
$ws.Cells.Item(3, 10).Value = 0.41
$ws.Cells.Item(3, 11).Value = 15
$ws.Cells.Item(3, 12).Value = 150
$ws.Cells.Item(3, 13).Value = 10
$ws.Cells.Item(3, 14).Value = 10
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = 0
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 0
$ws.Cells.Item(3, 19).Value = 0
$ws.Rows.Item(3).AutoFit()

# Row 4
$ws.Cells.Item(4, 1).Value = 45791.59371015046
$ws.Cells.Item(4, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4, 2).Value = $descNew
$ws.Cells.Item(4, 3).Value = $codeNew1
$ws.Cells.Item(4, 5).Value = $mockImproved
$ws.Cells.Item(4, 6).Value = 10
$ws.Cells.Item(4, 7).Value = 10
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 10
$ws.Cells.Item(4, 14).Value = 10
$ws.Cells.Item(4, 15).Value = 0
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = 0
$ws.Cells.Item(4, 19).Value = 0
$ws.Rows.Item(4).AutoFit()

# Row 5
$ws.Cells.Item(5, 1).Value = 45791.60139868056
$ws.Cells.Item(5, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5, 2).Value = $descNew
$ws.Cells.Item(5, 3).Value = $codeNew2
$ws.Cells.Item(5, 5).Value = $mockImproved
$ws.Cells.Item(5, 6).Value = 10
$ws.Cells.Item(5, 7).Value = 9
$ws.Cells.Item(5, 8).Value = -0.35
$ws.Cells.Item(5, 9).Value = -15.87
$ws.Cells.Item(5, 10).Value = 0.15
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 20
$ws.Cells.Item(5, 13).Value = 10
$ws.Cells.Item(5, 14).Value = 10
$ws.Cells.Item(5, 15).Value = 0
$ws.Cells.Item(5, 16).Value = 0
$ws.Cells.Item(5, 17).Value = 0
$ws.Cells.Item(5, 18).Value = 0
$ws.Cells.Item(5, 19).Value = 0
$ws.Rows.Item(5).AutoFit()

# Row 6
$ws.Cells.Item(6, 1).Value = 45791.60420985582
$ws.Cells.Item(6, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(6, 2).Value = $descRsi
$ws.Cells.Item(6, 3).Value = $codeRsi
$ws.Cells.Item(6, 5).Value = $mockImproved
$ws.Cells.Item(6, 6).Value = 10
$ws.Cells.Item(6, 7).Value = 10
$ws.Cells.Item(6, 8).Value = 2.73
$ws.Cells.Item(6, 9).Value = 49.95
$ws.Cells.Item(6, 10).Value = 0.41
$ws.Cells.Item(6, 11).Value = 15
$ws.Cells.Item(6, 12).Value = 150
$ws.Cells.Item(6, 13).Value = 10
$ws.Cells.Item(6, 14).Value = 10
$ws.Cells.Item(6, 15).Value = 0
$ws.Cells.Item(6, 16).Value = 0
$ws.Cells.Item(6, 17).Value = 0
$ws.Cells.Item(6, 18).Value = 0
$ws.Cells.Item(6, 19).Value = 0
$ws.Rows.Item(6).AutoFit()

